# edit.ps1
# Applies the OSTEOPOROSE_PATCH0_5 dashboard update to DASHBOARD.xlsx
# Sheets: "Slides" (sheet1) and "Batches" (sheet2)

$wb = $excel.ActiveWorkbook
$wsSlides  = $wb.Worksheets.Item("Slides")
$wsBatches = $wb.Worksheets.Item("Batches")

# ---------------------------------------------------------------------------
# Helper: write a plain-text value into a cell while (a) never letting the
# host auto-convert a date-looking string ("2026-01-25") into a date serial,
# and (b) leaving the cell's existing formatting/style completely untouched.
# We force the cell to "Text" number format before writing, then immediately
# restore the original formatting by pasting formats copied from a sibling
# cell that already carries the desired style (border/alignment/etc.).
# NOTE: this runtime's PowerShell subset does not bind named (-Param value)
# arguments on user-defined functions, so every call below uses positional
# args; pass $null explicitly for StyleSourceCell when there is none.
# ---------------------------------------------------------------------------
function Set-TextValue($Cell, $Value, $StyleSourceCell) {
    $Cell.NumberFormat = "@"
    $Cell.Value = $Value
    if ($StyleSourceCell -ne $null) {
        $StyleSourceCell.Copy()
        $Cell.PasteSpecial(-4122)  # xlPasteFormats
    }
}

# ===========================================================================
# Sheet "Slides": row 2 edits
# ===========================================================================
$wsSlides.Cells.Item(2,6).Value = "P1"
$wsSlides.Cells.Item(2,8).Value = "Figura incluída (KM esquemático) + citação SCOT-HEART 10y corrigida; padding ajustado."
$wsSlides.Cells.Item(2,9).Value = "Lancet 2025 (SCOT-HEART 10y; 10.1016/S0140-6736(24)01899-5); Circulation 2020 (LAP)"

# ===========================================================================
# Sheet "Slides": row 7 edits
# ===========================================================================
$wsSlides.Cells.Item(7,5).Value = "Reposicionado"
# G7 holds a date-looking string; keep it as text with its original style (s=2),
# using A7 (same row, style s=2) as the formatting donor.
Set-TextValue $wsSlides.Cells.Item(7,7) "2026-01-23" $wsSlides.Cells.Item(7,1)
$wsSlides.Cells.Item(7,8).Value = "Mantido conteúdo; reposicionado como ponte para EtD"
$wsSlides.Cells.Item(7,9).Value = "NEJM 2023 (CLEAR Outcomes)"

# ===========================================================================
# Sheet "Slides": rows 30 & 31 are fully replaced (new content + styling is
# cleared to the workbook default "Normal" style - no s= attribute at all).
# ===========================================================================

# --- Row 30 ---
$wsSlides.Cells.Item(30,1).Value = "OSTEOPOROSE"
$wsSlides.Cells.Item(30,2).Value = "VIEWER"
$wsSlides.Cells.Item(30,3).Value = "Viewer: fit-to-screen + safe bottom"
$wsSlides.Cells.Item(30,4).Value = "Infra"
$wsSlides.Cells.Item(30,5).Value = "Atualizado"
$wsSlides.Cells.Item(30,6).Value = "P0"
Set-TextValue $wsSlides.Cells.Item(30,7) "2026-01-26" $null
$wsSlides.Cells.Item(30,8).Value = "scheduleFit multi-pass + assets listeners + 100dvh"
$wsSlides.Cells.Item(30,9).Value = "—"
$wsSlides.Cells.Item(30,10).Value = "Validar em fullscreen/projetor"
$wsSlides.Range("A30:J30").Style = "Normal"

# --- Row 31 ---
$wsSlides.Cells.Item(31,1).Value = "OSTEOPOROSE"
$wsSlides.Cells.Item(31,2).Value = "PRINT"
$wsSlides.Cells.Item(31,3).Value = "Print/PDF 16:9: sizing + fit"
$wsSlides.Cells.Item(31,4).Value = "Infra"
$wsSlides.Cells.Item(31,5).Value = "Atualizado"
$wsSlides.Cells.Item(31,6).Value = "P0"
Set-TextValue $wsSlides.Cells.Item(31,7) "2026-01-26" $null
$wsSlides.Cells.Item(31,8).Value = "min-height override + print-fit.js + print.html regenerado"
$wsSlides.Cells.Item(31,9).Value = "—"
$wsSlides.Cells.Item(31,10).Value = "Re-testar export no Chrome/Edge"
$wsSlides.Range("A31:J31").Style = "Normal"

# ===========================================================================
# Sheet "Slides": brand-new rows 32 & 33 (default "Normal" style, no s=)
# ===========================================================================

# --- Row 32 ---
$wsSlides.Cells.Item(32,1).Value = "OSTEOPOROSE"
$wsSlides.Cells.Item(32,2).Value = "S08"
$wsSlides.Cells.Item(32,3).Value = "O que é Utilidade em Saúde?"
$wsSlides.Cells.Item(32,4).Value = "Utilidade/QALY"
$wsSlides.Cells.Item(32,5).Value = "Atualizado"
$wsSlides.Cells.Item(32,6).Value = "P1"
Set-TextValue $wsSlides.Cells.Item(32,7) "2026-01-26" $null
$wsSlides.Cells.Item(32,8).Value = "Redução de densidade + headings em navy (gold como acento)"
$wsSlides.Cells.Item(32,9).Value = "Brazier 2002; Peasgood 2009"
$wsSlides.Cells.Item(32,10).Value = "Checar se ainda precisa de fit"
$wsSlides.Range("A32:J32").Style = "Normal"

# --- Row 33 ---
$wsSlides.Cells.Item(33,1).Value = "OSTEOPOROSE"
$wsSlides.Cells.Item(33,2).Value = "S14"
$wsSlides.Cells.Item(33,3).Value = "Paradoxo: fraturas em osteopenia"
$wsSlides.Cells.Item(33,4).Value = "FRAX / Conceitos"
$wsSlides.Cells.Item(33,5).Value = "Atualizado"
$wsSlides.Cells.Item(33,6).Value = "P1"
Set-TextValue $wsSlides.Cells.Item(33,7) "2026-01-26" $null
$wsSlides.Cells.Item(33,8).Value = "Redução de densidade (padding/margens) + gráfico menor"
$wsSlides.Cells.Item(33,9).Value = "Siris 2004"
$wsSlides.Cells.Item(33,10).Value = "Checar legibilidade do gráfico em projeção"
$wsSlides.Range("A33:J33").Style = "Normal"

# ===========================================================================
# Sheet "Batches": brand-new row 5 (default "Normal" style, no s=)
# ===========================================================================
Set-TextValue $wsBatches.Cells.Item(5,1) "2026-01-26" $null
$wsBatches.Cells.Item(5,2).Value = "OSTEOPOROSE_PATCH0_5"
$wsBatches.Cells.Item(5,3).Value = "Viewer fit (corte inferior) + Print/PDF 16:9 (min-height) + polish S08/S14"
$wsBatches.Cells.Item(5,4).Value = "OSTEOPOROSE-changelog-dashboard_PATCH0_5.zip"
$wsBatches.Cells.Item(5,5).Value = "P0: stage 100dvh + safe bottom; print-fit.js; print.html regenerado"
$wsBatches.Range("A5:E5").Style = "Normal"

Write-Output "edit.ps1 completed"
